$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.069490432739258
$ws.Range("B1").Value = 6.327652454376221
$ws.Range("C1").Value = 6.441517353057861
$ws.Range("D1").Value = 6.858065128326416
$ws.Range("E1").Value = 5.033783435821533
